# qa_template_map.xlsx update
# - Add 5 new field-mapping rows to the "studies" (documents/studies/subjects/series)
#   mapping table for aerosol particle qc_validation fields.
# - Update the sheet selection/scroll position to reflect the new bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended after the existing last row (113) -> rows 114-118.
# Column A is always the table/group name ("studies"); columns B and C repeat
# the destination field name (same mapping "from" -> "to").
$newRows = @(
    @("studies", "aerosol_particle_diameter_mean"),
    @("studies", "aerosol_particle_diameter_gsd"),
    @("studies", "aerosol_particle_diameter_units"),
    @("studies", "aerosol_particle_density"),
    @("studies", "aerosol_particle_density_units")
)

$startRow = 114
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $group = $newRows[$i][0]
    $field = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $group
    $ws.Cells.Item($r, 2).Value = $field
    $ws.Cells.Item($r, 3).Value = $field
}

# Move the visible selection to the newly added second row (A115), matching
# the post-edit cursor position recorded in the workbook.
$ws.Range("A115").Select()
